# Apply weekly price update: insert two new daily records (row 138 and 139)
# above the existing data block, pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 138; this shifts rows
# 138:195 down to 140:197 and extends the used range / dimension
# automatically (A1:R195 -> A1:R197).
$ws.Rows("138:139").Insert()

# ---- Row 138 (new) ----
$ws.Range("A138").Value = 1
$ws.Range("B138").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C138").Value = "Arica y Parinacota"
$ws.Range("D138").Value = 44489
$ws.Range("E138").Value = 15
$ws.Range("F138").Value = 100112043
$ws.Range("G138").Value = "Pepino ensalada"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 120
$ws.Range("K138").Value = 5000
$ws.Range("L138").Value = 6000
$ws.Range("M138").Value = 5500
$ws.Range("N138").Value = "$/caja 70 unidades"
$ws.Range("O138").Value = "Región de Arica y Parinacota"
$ws.Range("P138").Value = 79
$ws.Range("Q138").Value = 70
$ws.Range("R138").Value = "Hortaliza"

# ---- Row 139 (new) ----
$ws.Range("A139").Value = 1
$ws.Range("B139").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C139").Value = "Arica y Parinacota"
$ws.Range("D139").Value = 44489
$ws.Range("E139").Value = 15
$ws.Range("F139").Value = 100112043
$ws.Range("G139").Value = "Pepino ensalada"
$ws.Range("H139").Value = "Sin especificar"
$ws.Range("I139").Value = "Segunda"
$ws.Range("J139").Value = 140
$ws.Range("K139").Value = 4000
$ws.Range("L139").Value = 5000
$ws.Range("M139").Value = 4500
$ws.Range("N139").Value = "$/caja 100 unidades"
$ws.Range("O139").Value = "Región de Arica y Parinacota"
$ws.Range("P139").Value = 45
$ws.Range("Q139").Value = 100
$ws.Range("R139").Value = "Hortaliza"

# Make sure the date cells keep the custom date number format used by
# the rest of the column (style carried over automatically from the
# insert, but set explicitly to be safe).
$ws.Range("D138:D139").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Done"
